$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay text (matches source inlineStr values)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "62.816.85"
$ws.Range("E2").Value = "  +0.20%  "

$ws.Range("D3").Value = "2.460.05"
$ws.Range("E3").Value = "  +0.63%  "

$ws.Range("D5").Value = "573.19"
$ws.Range("E5").Value = "  -0.69%  "

$ws.Range("D6").Value = "146.16"
$ws.Range("E6").Value = "  +0.36%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  -0.27%  "

$ws.Range("D9").Value = "2.460.43"
$ws.Range("E9").Value = "  +0.72%  "

$ws.Range("E10").Value = "  +0.54%  "

$ws.Range("D11").Value = "0.162"
$ws.Range("E11").Value = "  +0.22%  "

$ws.Range("E12").Value = "  +1.14%  "

$ws.Range("E13").Value = "  +1.51%  "

$ws.Range("D14").Value = "28.84"
$ws.Range("E14").Value = "  +2.35%  "

$ws.Range("D15").Value = "0.0000177"
$ws.Range("E15").Value = "  -0.66%  "

$ws.Range("D16").Value = "2.906.54"
$ws.Range("E16").Value = "  +0.56%  "

$ws.Range("D17").Value = "62.748.42"
$ws.Range("E17").Value = "  +0.39%  "

$ws.Range("D18").Value = "2.463.08"
$ws.Range("E18").Value = "  +0.54%  "

$ws.Range("E19").Value = "  +1.56%  "

$ws.Range("D20").Value = "11.00"
$ws.Range("E20").Value = "  +0.67%  "

$ws.Range("D21").Value = "326.84"
$ws.Range("E21").Value = "  -0.36%  "

$ws.Range("D22").Value = "2.22"
$ws.Range("E22").Value = "  +10.22%  "

$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("E24").Value = "  -0.03%  "

$ws.Range("D25").Value = "10.21"
$ws.Range("E25").Value = "  +20.48%  "

$ws.Range("D26").Value = "65.59"
$ws.Range("E26").Value = "  -0.16%  "

$ws.Range("D27").Value = "656.76"
$ws.Range("E27").Value = "  +2.55%  "

$ws.Range("D28").Value = "0.0₃0978"
$ws.Range("E28").Value = "  -0.23%  "

$ws.Range("E29").Value = "  +0.54%  "

$ws.Range("E30").Value = "  -14.39%  "

$ws.Range("D31").Value = "1.44"
$ws.Range("E31").Value = "  +0.34%  "

$ws.Range("D32").Value = "8.00"
$ws.Range("E32").Value = "  -2.18%  "

$ws.Range("E33").Value = "  -1.57%  "

$ws.Range("E34").Value = "  -2.65%  "

$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.04%  "

$ws.Range("D36").Value = "1.54"
$ws.Range("E36").Value = "  +3.30%  "

$ws.Range("E37").Value = "  +0.05%  "

$ws.Range("D38").Value = "5.40"
$ws.Range("E38").Value = "  -1.04%  "

$ws.Range("D39").Value = "0.368"
$ws.Range("E39").Value = "  -1.20%  "

$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").Value = "151.42"
$ws.Range("E40").Value = "  -0.50%  "

$ws.Range("B41").Value = "EthereumClassic"
$ws.Range("C41").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D41").Value = "18.69"
$ws.Range("E41").Value = "  +0.07%  "

$ws.Range("D42").Value = "2.74"
$ws.Range("E42").Value = "  +2.53%  "

$ws.Range("D43").Value = "1.73"
$ws.Range("E43").Value = "  -1.02%  "

$ws.Range("E44").Value = "  -74.06%  "

$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("D46").Value = "152.67"
$ws.Range("E46").Value = "  +5.08%  "

$ws.Range("E47").Value = "  +1.63%  "

$ws.Range("D48").Value = "3.57"
$ws.Range("E48").Value = "  -0.62%  "

$ws.Range("D49").Value = "20.55"
$ws.Range("E49").Value = "  -0.09%  "

$ws.Range("D50").Value = "0.606"
$ws.Range("E50").Value = "  +0.49%  "

$ws.Range("D51").Value = "0.0512"
$ws.Range("E51").Value = "  -0.48%  "
